$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

function Replace-InCell($row, $col, $findText, $replaceText) {
    $t = $d.Tables(1)
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                       $true, 0, $false, $replaceText, 1) | Out-Null
}

# Header date
Replace-All "2023-09-15 Friday" "2023-09-16 Saturday"

# Row block 1 (table rows 1..1, cols 1..5)
Replace-All "57×99=" "60×49="
Replace-All "75×55=" "80×95="
Replace-All "41×43=" "49×55="
Replace-All "32×69=" "64×32="
Replace-All "26×38=" "99×12="

# Row block 2 (table row 5)
Replace-All "40×49=" "50×82="
Replace-All "21×51=" "69×96="
Replace-All "37×42=" "66×18="
Replace-All "21×83=" "90×83="
Replace-All "16×30=" "82×12="

# Row block 3 (table row 10)
Replace-All "51×94=" "41×15="
Replace-All "86×16=" "24×14="
Replace-All "26×22=" "74×46="
Replace-All "46×99=" "14×80="
Replace-All "44×43=" "71×85="

# Row block 4 (table row 15) - "68×67=" appears twice (cols 3 and 5),
# so those two must be targeted per-cell rather than with a global replace.
Replace-All "23×26=" "42×67="
Replace-All "54×89=" "34×54="
Replace-InCell 15 3 "68×67=" "13×19="
Replace-All "75×25=" "44×56="
Replace-InCell 15 5 "68×67=" "80×45="

# Row block 5 (table row 20)
Replace-All "46×39=" "70×35="
Replace-All "16×11=" "74×37="
Replace-All "27×22=" "22×25="
Replace-All "18×68=" "53×84="
Replace-All "62×98=" "40×61="

Write-Host "Edits applied"
